# Update the cryptocurrency price/volume snapshot on Sheet1.
# Values that look like plain numbers are written with a leading "'" so
# Excel keeps them as literal text (matching the source data, which stores
# prices/percentages as text, not numbers — some prices even contain
# multiple "." thousands separators that would not parse as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.226.13"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "3.286.73"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'214.15"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").Value = "'630.04"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.383"
$ws.Range("E7").Value = "  +20.11%  "
$ws.Range("E8").Value = "  +15.95%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "3.286.78"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.578"
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("E12").Value = "  +12.73%  "
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -3.30%  "
$ws.Range("D14").Value = "'34.30"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "3.887.40"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "89.080.13"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "3.293.74"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'3.13"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "'14.17"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'438.38"
$ws.Range("D22").Value = "'8.90"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "'5.38"
$ws.Range("E23").Value = "  +3.21%  "
$ws.Range("D24").Value = "'7.39"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'5.26"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'12.35"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("D27").Value = "3.450.00"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").Value = "'76.86"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "'0.0000135"
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'0.188"
$ws.Range("E31").Value = "  +11.68%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'8.89"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("D34").Value = "'567.74"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("E35").Value = "  -9.15%  "
$ws.Range("D36").Value = "'7.19"
$ws.Range("E36").Value = "  +11.64%  "
$ws.Range("D37").Value = "'1.97"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = "  -7.44%  "
$ws.Range("D39").Value = "'22.72"
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'21.84"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").Value = "'3.09"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'0.401"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").Value = "'2.04"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'154.40"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "'181.40"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").Value = "'45.11"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'1.31"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.127"
$ws.Range("E50").Value = "  +15.08%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").Value = "'0.0681"
$ws.Range("E51").Value = "  +22.11%  "
